# chore: adapt column header formatting to respective input file names
#
# Renames the "_old" / "_new" header-column suffixes to the format-version
# specific suffixes "_FV2210" / "_FV2304", wraps the data range in an Excel
# Table ("Table1"), and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) -----------------------------------------
# Columns A:J carried the "_old" suffix -> "_FV2210"
# Columns L:U carried the "_new" suffix -> "_FV2304"
# Column K ("diff") is left untouched.

$newHeadersFV2210 = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)

$newHeadersFV2304 = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt 10; $i++) {
    # columns 1..10 => A..J
    $ws.Cells.Item(1, $i + 1).Value = $newHeadersFV2210[$i]
    # columns 12..21 => L..U
    $ws.Cells.Item(1, $i + 12).Value = $newHeadersFV2304[$i]
}

# --- 2. Freeze the header row ----------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Wrap the data range in an Excel Table -------------------------------
$rng = $ws.Range("A1:U62")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
